$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "246.60"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "0.98%"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "29.77"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "9.67%"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "1.28%"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.05714"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "0.95%"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "6.611"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "2.15%"
$c = $ws.Range("B7")
$c.NumberFormat = "@"
$c.Value = "MXToken"
$c = $ws.Range("C7")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.8581"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "4.52%"
$c = $ws.Range("B8")
$c.NumberFormat = "@"
$c.Value = "FTXToken"
$c = $ws.Range("C8")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.8671"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "2.89%"
$c = $ws.Range("B9")
$c.NumberFormat = "@"
$c.Value = "WazirX"
$c = $ws.Range("C9")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.1365"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "2.94%"
$c = $ws.Range("B10")
$c.NumberFormat = "@"
$c.Value = "MandalaExchangeToken"
$c = $ws.Range("C10")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07068"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "2.15%"
$c = $ws.Range("B11")
$c.NumberFormat = "@"
$c.Value = "BitrueCoin"
$c = $ws.Range("C11")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.02861"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "-4.23%"
$c = $ws.Range("B12")
$c.NumberFormat = "@"
$c.Value = "BitMartToken"
$c = $ws.Range("C12")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.09390"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "0.10%"
$c = $ws.Range("B13")
$c.NumberFormat = "@"
$c.Value = "BitForexToken"
$c = $ws.Range("C13")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.001534"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "1.27%"
$c = $ws.Range("B14")
$c.NumberFormat = "@"
$c.Value = "CoinExToken"
$c = $ws.Range("C14")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.04157"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "-1.24%"
$c = $ws.Range("B15")
$c.NumberFormat = "@"
$c.Value = "One"
$c = $ws.Range("C15")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0005981"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "-0.17%"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.006139"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "-0.30%"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "3,764.37%"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "-0.96%"
$c = $ws.Range("B19")
$c.NumberFormat = "@"
$c.Value = "GateToken"
$c = $ws.Range("C19")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.051"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "1.65%"
$c = $ws.Range("B20")
$c.NumberFormat = "@"
$c.Value = "BTSEToken"
$c = $ws.Range("C20")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "2.182"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "-1.88%"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.3145"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "1.03%"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.03250"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "2.91%"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.1300"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "3.66%"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "3.477"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "-2.24%"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.005092"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "14.22%"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.001222"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "-0.08%"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.0001210"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "23.50%"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.03768"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "2.68%"
$c = $ws.Range("B41")
$c.NumberFormat = "@"
$c.Value = "BKEXToken"
$c = $ws.Range("C41")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.1070"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "1.58%"
$c = $ws.Range("B42")
$c.NumberFormat = "@"
$c.Value = "CEJI"
$c = $ws.Range("C42")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.002531"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "-0.37%"
$c = $ws.Range("B43")
$c.NumberFormat = "@"
$c.Value = "KickToken"
$c = $ws.Range("C43")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.003481"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "-42.46%"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.009402"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "13.45%"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.00005108"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "-3.96%"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "0.00%"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.07512"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "-41.42%"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.002724"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "5.62%"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.00002100"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "0.00%"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0002000"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "0.00%"
